$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Values for columns G6..G9 (table columns 7..10) across the 5 data rows (table rows 2..6)
$values = @(
    @(7, 8, 10, 9),
    @(5, 7, 10, 10),
    @(7, 7, 10, 9),
    @(7, 8, 9, 8),
    @(6, 6, 8, 8)
)

for ($r = 0; $r -lt 5; $r++) {
    $tableRow = $r + 2
    for ($c = 0; $c -lt 4; $c++) {
        $tableCol = $c + 7
        $cell = $t.Cell($tableRow, $tableCol)

        if ($r -eq 4 -and $c -eq 0) {
            # Last edited cell (G6 / Oral expression row) keeps the editing
            # cursor mark ("_GoBack") that Word drops at the most recent edit.
            $rng = $cell.Range
            $rng.Collapse(0)
            $d.Bookmarks.Add("_GoBack", $rng)
        }

        $cell.Range.Text = [string]$values[$r][$c]
    }
}
